# Scorecard Knowledge: add shortcut name
# - Fix minor wording in existing full-name columns (B,C) for module 1 & 3
# - Add two new columns: shortcut_name_km (D) and shortcut_name_en (E)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScorecardKnowledge")

# --- Minor text corrections on existing columns ---
$ws.Range("B2").Value = "មេរៀនម៉ូឌុលទី ១៖ ការណែនាំអំពីគណនេយ្យភាពសង្គម (ISAF)"
$ws.Range("C2").Value = "Module 1: introduction to ISAF"
$ws.Range("C4").Value = "Module 3: Facilitating community scorecard and service provider self-assessment"

# --- New header cells ---
$ws.Range("D1").Value = "shortcut_name_km"
$ws.Range("E1").Value = "shortcut_name_en"
$ws.Range("D1").Font.Bold = $true
$ws.Range("E1").Font.Bold = $true

# --- New shortcut-name data rows ---
$ws.Range("D2").Value = "មេរៀនម៉ូឌុលទី ១"
$ws.Range("E2").Value = "Module 1"

$ws.Range("D3").Value = "មេរៀនម៉ូឌុលទី ២"
$ws.Range("E3").Value = "Module 2"

$ws.Range("D4").Value = "មេរៀនម៉ូឌុលទី ៣"
$ws.Range("E4").Value = "Module 3"

$ws.Range("D5").Value = "មេរៀនម៉ូឌុលទី ៤"
$ws.Range("E5").Value = "Module 4"

# --- Column widths to fit the new content ---
$ws.Columns.Item(2).ColumnWidth = 41.6015625
$ws.Columns.Item(3).ColumnWidth = 40.33203125
$ws.Columns.Item(4).ColumnWidth = 15.33203125
